$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "[-, -, -, 'MCT-3A-Robótica']"
$ws.Range("E2").Value = "[-, -, 'MCT-3A-Robótica', -]"
$ws.Range("E3").Value = "-"
$ws.Range("B6").Value = "-"
$ws.Range("B7").Value = "[-, 'MCT-3A-Robótica', -, -]"
$ws.Range("E8").Value = "[-, 'MCT-3A-Robótica', -, -]"
$ws.Range("E18").Value = "-"
$ws.Range("F20").Value = "['ELM-2NA-Eletrônica Básica', -]"
$ws.Range("E21").Value = "-"
$ws.Range("F21").Value = "['ELM-2NA-Eletrônica Básica', -]"
